# Add a new library-formula entry: "checkAnalysisUnitStatus"
# Sheet "Library_Formula" gains row 16 (Action=CREATE/MODIFY, Library=LIB_EWS,
# Formula Name=checkAnalysisUnitStatus, Output=String, Input=String).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row 16 so it inherits row 15's per-column formatting
# (same trick Excel itself uses when a row is appended under the last one).
$ws.Rows("16").Insert() | Out-Null

$ws.Range("A16").Value = "CREATE/MODIFY"
$ws.Range("B16").Value = "LIB_EWS"
$ws.Range("C16").Value = "checkAnalysisUnitStatus"
$ws.Range("E16").Value = "String"
$ws.Range("F16").Value = "String"

# C15's inherited style differs from the rest of the column (it uses the
# "bestFit" xf); match the plain Trebuchet MS 10pt style used elsewhere.
$ws.Range("C16").Font.Name = "Trebuchet MS"
$ws.Range("C16").Font.Size = 10

# Selection moves on one row, same as it did for every previous entry.
$ws.Range("C19").Select() | Out-Null
